# Atualização automática (19/02/2026  8:15:06,85)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 63 (pushes rows 63..77 down to 64..78),
# copying formatting from the row above (Excel default insert behaviour).
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new rejected-item record.
$ws.Cells.Item(63, 2).Value = "180- CX-360G"
$ws.Cells.Item(63, 3).Value = "18/02/2026"
$ws.Cells.Item(63, 4).Value = (Get-Date -Year 2026 -Month 2 -Day 18 -Hour 14 -Minute 26 -Second 5)
$ws.Cells.Item(63, 5).Value = (Get-Date -Year 2026 -Month 2 -Day 18 -Hour 14 -Minute 26 -Second 5)
$ws.Cells.Item(63, 6).Value = "14"
$ws.Cells.Item(63, 7).Value = "2"
$ws.Cells.Item(63, 8).Value = "0505 - FALHA NO SLEEVE"
$ws.Cells.Item(63, 9).Value = 1418251
$ws.Cells.Item(63, 10).Value = "CIPF002225-MINNIE_26-24 - GARRAFA RETRÔ 1L SOPRO - IMPRESSÃO DIGIT"
$ws.Cells.Item(63, 11).Value = "8064 - MATHEUS ANZOLIN"
$ws.Cells.Item(63, 12).Value = "24"
$ws.Cells.Item(63, 14).Value = 1

# Update the trailing summary rows (now shifted to 77 and 78) to reflect
# the extra rejected-record count.
$ws.Cells.Item(77, 14).Value = 144
$ws.Cells.Item(77, 15).Value = 30
$ws.Cells.Item(78, 8).Value = 71
$ws.Cells.Item(78, 14).Value = 316
$ws.Cells.Item(78, 15).Value = 71
